$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.718.83'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -4.39%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.646.45'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -6.50%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.63%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.26'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.94%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3635'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -5.46%  '

$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.71'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -8.85%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3282'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -10.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.129'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -8.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07054'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -7.97%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.006'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.023'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -7.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.67'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -9.97%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.651'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -6.91%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.645.05'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -6.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001061'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -8.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06590'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.83%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.73%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '79.39'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -9.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.24'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -8.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.998'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -8.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.15'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -5.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.660.84'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.56%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.415'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.542'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -14.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.34'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.19'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -7.75%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '127.43'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -6.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.832.05'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -6.70%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.102'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -8.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.123'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.40%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.090'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -15.63%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.725'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.53%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08447'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.61'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -10.25%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.218'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -8.37%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06195'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -8.71%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02279'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -8.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2089'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -6.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.223'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -6.51%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.314'
$ws.Range("D42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6063'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -8.17%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.005'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.86%  '

$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.758'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.22%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.01'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -8.83%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5745'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -9.97%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.20'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -7.82%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.978'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -9.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07037'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -6.20%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.77'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -7.79%  '
